# edit.ps1 -- apply "New crime data collected" update to CompStat_1 sheet
# Updates the report header (volume/number + week-range) and the weekly
# crime-complaint figures (Week to Date / 28 Day / Year to Date columns,
# plus their derived % Chg columns) for rows 16-28.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Report header: issue number and week-covering date range ---
$ws.Range("A8").Value = "Volume 31   Number  24"
$ws.Range("C9").Value = "Report Covering the Week  6/10/2024  Through  6/16/2024"

# --- Row 16 (Robbery) ---
$ws.Range("C16").Value = 1
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 3
$ws.Range("H16").Value = 66.666666666666
$ws.Range("I16").Value = 36
$ws.Range("K16").Value = 9.090909090909
$ws.Range("L16").Value = 5.882352941176
$ws.Range("M16").Value = -20
$ws.Range("N16").Value = -88.888888888888

# --- Row 17 (Fel. Assault) ---
$ws.Range("C17").Value = 3
$ws.Range("E17").Value = 200
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = -11.111111111111
$ws.Range("I17").Value = 44
$ws.Range("J17").Value = 46
$ws.Range("K17").Value = -4.347826086956
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 57.142857142857
$ws.Range("N17").Value = -21.428571428571

# --- Row 18 (Burglary) ---
$ws.Range("C18").Value = 3
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 20
$ws.Range("I18").Value = 43
$ws.Range("K18").Value = -21.818181818181
$ws.Range("L18").Value = -14
$ws.Range("M18").Value = -27.118644067796
$ws.Range("N18").Value = -92.833333333333

# --- Row 19 (Gr. Larceny) ---
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 100
$ws.Range("G19").Value = 41
$ws.Range("H19").Value = -7.317073170731
$ws.Range("I19").Value = 193
$ws.Range("J19").Value = 222
$ws.Range("K19").Value = -13.063063063063
$ws.Range("L19").Value = -20.247933884297
$ws.Range("M19").Value = 22.151898734177
$ws.Range("N19").Value = -56.919642857142

# --- Row 20 (G.L.A.) ---
$ws.Range("C20").Value = 4
$ws.Range("E20").Value = 33.333333333333
$ws.Range("F20").Value = 26
$ws.Range("H20").Value = 85.714285714285
$ws.Range("I20").Value = 79
$ws.Range("J20").Value = 56
$ws.Range("K20").Value = 41.071428571428
$ws.Range("L20").Value = 113.513513513514
$ws.Range("M20").Value = 58
$ws.Range("N20").Value = -95.018915510718

# --- Row 21 (TOTAL) ---
$ws.Range("C21").Value = 23
$ws.Range("D21").Value = 10
$ws.Range("E21").Value = 130
$ws.Range("F21").Value = 84
$ws.Range("G21").Value = 72
$ws.Range("H21").Value = 16.666666666666
$ws.Range("I21").Value = 401
$ws.Range("J21").Value = 414
$ws.Range("K21").Value = -3.140096618357
$ws.Range("L21").Value = -3.373493975903
$ws.Range("M21").Value = 17.595307917888
$ws.Range("N21").Value = -86.713055003313

# --- Row 22 (Transit) ---
$ws.Range("D22").Value = 1
$ws.Range("D22").NumberFormat = '#,##0'
$ws.Range("E22").Value = 100
$ws.Range("E22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F22").Value = 6
$ws.Range("G22").Value = 1
$ws.Range("G22").NumberFormat = '#,##0'
$ws.Range("H22").Value = 500
$ws.Range("H22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("I22").Value = 20
$ws.Range("J22").Value = 12
$ws.Range("K22").Value = 66.666666666666
$ws.Range("L22").Value = 33.333333333333
$ws.Range("M22").Value = 150

# --- Row 24 (Petit Larceny) ---
$ws.Range("C24").Value = 38
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = 26.666666666666
$ws.Range("F24").Value = 136
$ws.Range("G24").Value = 121
$ws.Range("H24").Value = 12.396694214876
$ws.Range("I24").Value = 739
$ws.Range("J24").Value = 728
$ws.Range("K24").Value = 1.510989010989
$ws.Range("L24").Value = -9.213759213759
$ws.Range("M24").Value = 72.261072261072

# --- Row 25 (Retail Theft) ---
$ws.Range("C25").Value = 24
$ws.Range("D25").Value = 26
$ws.Range("E25").Value = -7.692307692307
$ws.Range("F25").Value = 88
$ws.Range("G25").Value = 91
$ws.Range("H25").Value = -3.296703296703
$ws.Range("I25").Value = 538
$ws.Range("J25").Value = 527
$ws.Range("K25").Value = 2.087286527514
$ws.Range("L25").Value = -8.658743633276

# --- Row 26 (Misd. Assault) ---
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = 150
$ws.Range("G26").Value = 20
$ws.Range("H26").Value = 40
$ws.Range("I26").Value = 130
$ws.Range("J26").Value = 111
$ws.Range("K26").Value = 17.117117117117
$ws.Range("L26").Value = 52.941176470588
$ws.Range("M26").Value = 30

# --- Row 28 (Other Sex Crimes) ---
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 1
$ws.Range("D28").NumberFormat = '#,##0'
$ws.Range("E28").Value = 100
$ws.Range("E28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F28").Value = 7
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 600
$ws.Range("I28").Value = 16
$ws.Range("J28").Value = 14
$ws.Range("K28").Value = 14.285714285714
$ws.Range("L28").Value = -11.111111111111
